# Apply the STM32 Hero Board GPIO worksheet updates:
#  - Row 6  (PA4):  rename PX_GPIO_DAC       -> PX_GPIO_DAC1
#  - Row 21 (PB2):  rename PX_GPIO_SPI2_CS_DF -> PX_GPIO_SPI2_CS_SF
#  - Row 33 (PB14): change Pull from PX_GPIO_PULL_NO -> PX_GPIO_PULL_DN
#  - Row 46 (PC10): rename PX_GPIO_UART3_TX  -> PX_GPIO_UART4_TX, AltFn PX_GPIO_AF_4 -> PX_GPIO_AF_6
#  - Row 47 (PC11): rename PX_GPIO_UART3_RX  -> PX_GPIO_UART4_RX, AltFn PX_GPIO_AF_4 -> PX_GPIO_AF_6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPIO")

$ws.Range("J33").Value = "PX_GPIO_PULL_DN"

$ws.Range("A46").Value = "PX_GPIO_UART4_TX"
$ws.Range("N46").Value = "PX_GPIO_AF_6"

$ws.Range("A47").Value = "PX_GPIO_UART4_RX"
$ws.Range("N47").Value = "PX_GPIO_AF_6"

$ws.Range("A6").Value = "PX_GPIO_DAC1"

$ws.Range("A21").Value = "PX_GPIO_SPI2_CS_SF"
